$wb = $excel.ActiveWorkbook

# New header row (row 1) shared across all sheets
$headers = @("eb","gb","hp","st","wi","ieh","chp","ac","ab_ct","ab_hp","cp_ct","cp_hp","ttes","btes","ites")

# New row-2 data values, keyed by sheet name
$data = @{
    "2025" = @(3906.399109145206, 0, 48353.76274462014, 0, 289724.0114301849, 9433.134471502228, 0, 2534.277928792104, 0, 0, 0, 0, 0, 2366.658982337573, 1995.412676509708)
    "2030" = @(6991.052031681918, 0, 197913.7502057619, 0, 289724.0114301849, 16452.51445364119, 0, 8194.52068131253, 0, 0, 0, 0, 0, 7540.299372506026, 6256.367679449893)
    "2035" = @(31236.29455387744, 0, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 12882.34160925225, 9261.533324927314)
    "2040" = @(31236.29455387744, 0, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 14040.19136169543, 9261.533324927314)
    "2045" = @(38906.8534480406, 193.0947398408091, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 16870.34251754129, 10093.83191070589)
    "2050" = @(38906.8534480406, 193.0947398408091, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 16870.34251754129, 10093.83191070589)
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($data.ContainsKey($name)) {
        for ($i = 0; $i -lt $headers.Length; $i++) {
            $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
        }
        $rowVals = $data[$name]
        for ($i = 0; $i -lt $rowVals.Length; $i++) {
            $ws.Cells.Item(2, $i + 1).Value = $rowVals[$i]
        }
    }
}
